# Apply the updated crypto price/volume figures scraped on
# Mon Jan 15 17:39:47 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.668.64"
$ws.Range("D3").Value = "2.532.01"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'317.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.62%  "
$ws.Range("D6").Value = "'95.23"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.81%  "
$ws.Range("D7").Value = "'0.574"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.08%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("D10").Value = "'36.17"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("D11").Value = "'0.0808"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "'7.55"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("D14").Value = "2.921.41"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "'15.49"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "2.520.44"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "'0.850"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "42.655.60"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'12.99"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "0.0₃0962"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").Value = "'70.17"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").Value = "'251.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").Value = "'26.63"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "'2.40"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.81%  "
$ws.Range("D29").Value = "'39.28"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("D30").Value = "'10.17"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").Value = "'6.08"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").Value = "'155.75"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'2.12"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "'19.19"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.42%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("D36").Value = "'0.0786"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("D40").Value = "'23.65"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "'2.32"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +11.42%  "
$ws.Range("D42").Value = "'3.81"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("E45").Value = "  -5.70%  "
$ws.Range("D46").Value = "2.016.45"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").Value = "'85.66"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").Value = "2.775.65"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "'74.29"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "'102.72"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.98%  "
